$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "datos actualizados" timestamp (row 1)
$ws.Range("A1").Value = "Datos actualizados a 12 de Octubre de 2020 a las 15:34"

# Row 4: Estados Unidos
$ws.Range("A4").Value = "Estados Unidos"
$ws.Range("B4").Value = 7993215
$ws.Range("C4").Value = 1217
$ws.Range("D4").Value = 5128497
$ws.Range("E4").Value = 2645010
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 13
$ws.Range("H4").Value = 219708

# Row 5: India
$ws.Range("A5").Value = "India"
$ws.Range("B5").Value = 7128268
$ws.Range("C5").Value = 8968
$ws.Range("D5").Value = 6149535
$ws.Range("E5").Value = 869448
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 101
$ws.Range("H5").Value = 109285

# Row 18: Irak
$ws.Range("A18").Value = "Irak"
$ws.Range("B18").Value = 405437
$ws.Range("C18").Value = 3107
$ws.Range("D18").Value = 340050
$ws.Range("E18").Value = 55475
$ws.Range("F18").Value = 0
$ws.Range("G18").Value = 60
$ws.Range("H18").Value = 9912

# Row 22: Arabia Saudita
$ws.Range("A22").Value = "Arabia Saudita"
$ws.Range("B22").Value = 339615
$ws.Range("C22").Value = 348
$ws.Range("D22").Value = 325839
$ws.Range("E22").Value = 8708
$ws.Range("F22").Value = 0
$ws.Range("G22").Value = 25
$ws.Range("H22").Value = 5068

# Row 30: Paises Bajos
$ws.Range("A30").Value = "Paises Bajos"
$ws.Range("B30").Value = 181498
$ws.Range("C30").Value = 6845
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("F30").Value = 0
$ws.Range("G30").Value = 12
$ws.Range("H30").Value = 6596

# Row 50: Portugal
$ws.Range("A50").Value = "Portugal"
$ws.Range("B50").Value = 87913
$ws.Range("C50").Value = 1249
$ws.Range("D50").Value = 53498
$ws.Range("E50").Value = 32321
$ws.Range("F50").Value = 0
$ws.Range("G50").Value = 14
$ws.Range("H50").Value = 2094

# Row 51: Costa Rica
$ws.Range("A51").Value = "Costa Rica"
$ws.Range("B51").Value = 87439
$ws.Range("C51").Value = 0
$ws.Range("D51").Value = 52669
$ws.Range("E51").Value = 33694
$ws.Range("F51").Value = 0
$ws.Range("G51").Value = 0
$ws.Range("H51").Value = 1076

# Row 57: Barein
$ws.Range("A57").Value = "Barein"
$ws.Range("B57").Value = 75614
$ws.Range("C57").Value = 0
$ws.Range("D57").Value = 71249
$ws.Range("E57").Value = 4087
$ws.Range("F57").Value = 0
$ws.Range("G57").Value = 3
$ws.Range("H57").Value = 278

# Row 60: Uzbekistan
$ws.Range("A60").Value = "Uzbekistan"
$ws.Range("B60").Value = 61319
$ws.Range("C60").Value = 221
$ws.Range("D60").Value = 58306
$ws.Range("E60").Value = 2504
$ws.Range("F60").Value = 0
$ws.Range("G60").Value = 4
$ws.Range("H60").Value = 509

# Row 70: Estado de Palestina
$ws.Range("A70").Value = "Estado de Palestina"
$ws.Range("B70").Value = 44684
$ws.Range("C70").Value = 385
$ws.Range("D70").Value = 38228
$ws.Range("E70").Value = 6069
$ws.Range("F70").Value = 0
$ws.Range("G70").Value = 6
$ws.Range("H70").Value = 387

# Row 74: Kenia
$ws.Range("A74").Value = "Kenia"
$ws.Range("B74").Value = 41619
$ws.Range("C74").Value = 73
$ws.Range("D74").Value = 32000
$ws.Range("E74").Value = 8842
$ws.Range("F74").Value = 0
$ws.Range("G74").Value = 11
$ws.Range("H74").Value = 777

# Row 77: Serbia
$ws.Range("A77").Value = "Serbia"
$ws.Range("B77").Value = 34854
$ws.Range("C77").Value = 67
$ws.Range("D77").Value = 31536
$ws.Range("E77").Value = 2553
$ws.Range("F77").Value = 0
$ws.Range("G77").Value = 2
$ws.Range("H77").Value = 765

# Row 95: Noruega
$ws.Range("A95").Value = "Noruega"
$ws.Range("B95").Value = 15585
$ws.Range("C95").Value = 61
$ws.Range("D95").Value = 11863
$ws.Range("E95").Value = 3446
$ws.Range("F95").Value = 0
$ws.Range("G95").Value = 1
$ws.Range("H95").Value = 276

# Row 96: Zambia
$ws.Range("A96").Value = "Zambia"
$ws.Range("B96").Value = 15549
$ws.Range("C96").Value = 91
$ws.Range("D96").Value = 14682
$ws.Range("E96").Value = 522
$ws.Range("F96").Value = 0
$ws.Range("G96").Value = 8
$ws.Range("H96").Value = 345

# Row 107: Tayikistan
$ws.Range("A107").Value = "Tayikistan"
$ws.Range("B107").Value = 10260
$ws.Range("C107").Value = 38
$ws.Range("D107").Value = 9103
$ws.Range("E107").Value = 1078
$ws.Range("F107").Value = 0
$ws.Range("G107").Value = 0
$ws.Range("H107").Value = 79

# Row 110: Uganda
$ws.Range("A110").Value = "Uganda"
$ws.Range("B110").Value = 9864
$ws.Range("C110").Value = 63
$ws.Range("D110").Value = 6109
$ws.Range("E110").Value = 3661
$ws.Range("F110").Value = 0
$ws.Range("G110").Value = 1
$ws.Range("H110").Value = 94

# Row 139: Mayotte
$ws.Range("A139").Value = "Mayotte"
$ws.Range("B139").Value = 4030
$ws.Range("C139").Value = 41
$ws.Range("D139").Value = 2964
$ws.Range("E139").Value = 1023
$ws.Range("F139").Value = 0
$ws.Range("G139").Value = 0
$ws.Range("H139").Value = 43
